$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(15, 1, 27.142, 407.13, 2.14175312, 23150874.4, 6.699889079999999, 12.10993844, 341.42259628, 5.432205234286976, 81.48307851430465, 0.3962163935013594, 11838347.66101951, 1.183531674919023, 3.189758000172819, 158.7446962498986, 0.15),
    @(15, 1, 45.54, 683.014, 1.27201826, 31590122.696, 7.224497, 21.98348316, 1039.1950829, 8.657152683484931, 129.7456808043133, 0.2225436130582013, 12273163.91046628, 1.127609198688683, 5.610320518809024, 465.6026588659575, 0.85),
    @(20, 1, 19.82, 396.4, 2.18429414, 35823788.976, 8.88854418, 8.803746520000001, 180.07724766, 3.535052076896852, 70.70104153793703, 0.3629088317417651, 14446475.93031386, 1.64008603265614, 2.272922095012736, 78.81213722572571, 0.15),
    @(20, 1, 34.05, 680.926, 1.27813448, 53912738.104, 10.07716026, 17.1596059, 607.1171421399999, 6.679875858349909, 133.4753905548417, 0.2283522623697615, 19935099.88143703, 1.487449818530974, 4.289681515161005, 281.4396754879621, 0.85)
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
